# Update cryptos list: refresh price (D) and volume-1h (E) columns,
# and a handful of rows whose coin (B) / link (C) changed because the
# underlying ranking reshuffled position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to push numeric-looking price strings ("1.002",
# "274.10", ...) into the sheet as literal TEXT via Copy/PasteSpecial
# (values-only) instead of Range.Value, which would otherwise let Excel
# auto-coerce them into real numbers and touch their cell style.
$helper = $ws.Range("Z1")

$helper.Formula = "=""19.970.37"""
$helper.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -8.20%  "

$helper.Formula = "=""1.418.06"""
$helper.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -7.88%  "

$helper.Formula = "=""1.002"""
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.17%  "

$helper.Formula = "=""1.002"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.12%  "

$helper.Formula = "=""274.10"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -5.34%  "

$helper.Formula = "=""0.3709"""
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -4.78%  "

$helper.Formula = "=""0.3077"""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -3.58%  "

$helper.Formula = "=""39.54"""
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -8.20%  "

$helper.Formula = "=""1.004"""
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -5.13%  "

$helper.Formula = "=""0.06591"""
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -8.54%  "

$helper.Formula = "=""1.002"""
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.11%  "

$helper.Formula = "=""5.419"""
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -3.97%  "

$helper.Formula = "=""17.10"""
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -8.10%  "

$helper.Formula = "=""6.172"""
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -6.68%  "

$helper.Formula = "=""1.422.77"""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -7.66%  "

$helper.Formula = "=""0.00001007"""
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -9.60%  "

$helper.Formula = "=""0.05766"""
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -12.54%  "

$helper.Formula = "=""74.20"""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -10.79%  "

$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.12%  "

$helper.Formula = "=""5.618"""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -8.56%  "

$helper.Formula = "=""14.50"""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -5.91%  "

$helper.Formula = "=""10.99"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.77%  "

$helper.Formula = "=""2.333"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -3.40%  "

$helper.Formula = "=""19.987.19"""
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -8.15%  "

$helper.Formula = "=""2.285"""
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -3.89%  "

$helper.Formula = "=""139.18"""
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -5.04%  "

$helper.Formula = "=""16.95"""
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -7.68%  "

$helper.Formula = "=""1.582.02"""
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -7.75%  "

$helper.Formula = "=""109.21"""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -7.22%  "

$helper.Formula = "=""3.856"""
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -20.21%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$helper.Formula = "=""5.394"""
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -8.85%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$helper.Formula = "=""0.8767"""
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -9.72%  "

$helper.Formula = "=""0.07741"""
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -5.46%  "

$helper.Formula = "=""8.476"""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -4.49%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$helper.Formula = "=""0.05728"""
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -5.77%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$helper.Formula = "=""4.786"""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -6.99%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$helper.Formula = "=""10.84"""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +1.50%  "

$helper.Formula = "=""0.1933"""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -5.21%  "

$helper.Formula = "=""0.02048"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -7.24%  "

$helper.Formula = "=""1.070"""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -10.22%  "

$helper.Formula = "=""1.281"""
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -14.11%  "

$helper.Formula = "=""0.5317"""
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -7.55%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$helper.Formula = "=""12.28"""
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -5.84%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$helper.Formula = "=""3.532"""
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -5.62%  "

$helper.Formula = "=""0.5133"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -6.94%  "

$ws.Range("E48").Value = "  -3.63%  "

$helper.Formula = "=""109.20"""
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -7.07%  "

$helper.Formula = "=""1.049"""
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -10.01%  "

$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.06%  "

$helper.ClearContents()
$excel.CutCopyMode = 0
